$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text to lowercase to match new naming convention
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "name"

# Reflect the final cell selection/cursor position at save time
[void]$ws.Range("G8").Select()
